# Adds two new, empty "ListParagraph" paragraphs immediately after the
# paragraph that contains "console.log(details);" (and before the blank
# paragraph that already follows it), matching the formatting of the
# surrounding list paragraphs (Georgia 15pt, color 292929, white shading,
# no space-after).

$d = $word.ActiveDocument

# Locate the anchor paragraph's text and collapse the range to its end.
$rng = $d.Content
[void]$rng.Find.Execute("console.log(details);", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# OOXML for a single empty paragraph with the desired formatting. This is
# wrapped in the minimal pkg:package envelope that Range.InsertXML expects.
$emptyParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0"/><w:rPr><w:rFonts w:ascii="Georgia" w:eastAsia="Times New Roman" w:hAnsi="Georgia" w:cs="Times New Roman"/><w:color w:val="292929"/><w:spacing w:val="-1"/><w:sz w:val="30"/><w:szCs w:val="30"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Insert two new paragraph marks after "console.log(details);".
[void]$rng.InsertParagraphAfter()
$rng.Collapse(0)
[void]$rng.InsertParagraphAfter()

# The two InsertParagraphAfter calls above create well-formatted paragraph
# marks, but also leave a stray empty run carrying the inherited run
# properties. Replace each new paragraph's content with the clean OOXML
# fragment (no run) so the result matches a genuinely empty paragraph.
$n = $d.Paragraphs.Count
[void]$d.Paragraphs.Item($n - 2).Range.InsertXML($emptyParaXml)
[void]$d.Paragraphs.Item($n - 1).Range.InsertXML($emptyParaXml)
